# Update the cryptocurrency price table with refreshed values.
# Cells in column D hold price strings and column E holds percentage-change
# strings; both are formatted as text in the source data (they can contain
# multiple "." thousand separators, a leading/trailing space, "%" signs, or
# special Unicode digits), so we force text formatting before assigning the
# value to stop Excel from "helpfully" re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '37.775.45' },
    @{ Cell = 'E2'; Value = '  -1.16%  ' },
    @{ Cell = 'D3'; Value = '2.081.74' },
    @{ Cell = 'E3'; Value = '  -1.67%  ' },
    @{ Cell = 'E4'; Value = '  -0.04%  ' },
    @{ Cell = 'D5'; Value = '233.64' },
    @{ Cell = 'E5'; Value = '  -0.68%  ' },
    @{ Cell = 'D6'; Value = '0.625' },
    @{ Cell = 'E6'; Value = '  -0.31%  ' },
    @{ Cell = 'D7'; Value = '58.80' },
    @{ Cell = 'E7'; Value = '  +0.75%  ' },
    @{ Cell = 'E8'; Value = '  -0.06%  ' },
    @{ Cell = 'E9'; Value = '  +0.19%  ' },
    @{ Cell = 'E10'; Value = '  +0.60%  ' },
    @{ Cell = 'E11'; Value = '  +2.76%  ' },
    @{ Cell = 'B12'; Value = 'Chainlink' },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' },
    @{ Cell = 'D12'; Value = '14.87' },
    @{ Cell = 'E12'; Value = '  +0.79%  ' },
    @{ Cell = 'B13'; Value = 'Avalanche' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' },
    @{ Cell = 'D13'; Value = '21.09' },
    @{ Cell = 'E13'; Value = '  -3.05%  ' },
    @{ Cell = 'B14'; Value = 'Polygon' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' },
    @{ Cell = 'D14'; Value = '0.776' },
    @{ Cell = 'E14'; Value = '  -1.82%  ' },
    @{ Cell = 'B15'; Value = 'Polkadot' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Cell = 'D15'; Value = '5.36' },
    @{ Cell = 'E15'; Value = '  +1.47%  ' },
    @{ Cell = 'B16'; Value = 'WrappedEther' },
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Cell = 'D16'; Value = '2.073.88' },
    @{ Cell = 'E16'; Value = '  -2.59%  ' },
    @{ Cell = 'B17'; Value = 'WrappedBTC' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Cell = 'D17'; Value = '37.739.57' },
    @{ Cell = 'E17'; Value = '  -1.07%  ' },
    @{ Cell = 'B18'; Value = 'Uniswap' },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' },
    @{ Cell = 'D18'; Value = '6.13' },
    @{ Cell = 'E18'; Value = '  -0.63%  ' },
    @{ Cell = 'B19'; Value = 'Litecoin' },
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' },
    @{ Cell = 'D19'; Value = '71.66' },
    @{ Cell = 'E19'; Value = '  +1.09%  ' },
    @{ Cell = 'B20'; Value = 'ShibaInu' },
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Cell = 'D20'; Value = '0.0₃0836' },
    @{ Cell = 'E20'; Value = '  +1.03%  ' },
    @{ Cell = 'B21'; Value = 'BitcoinCash' },
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' },
    @{ Cell = 'D21'; Value = '228.55' },
    @{ Cell = 'E21'; Value = '  -0.24%  ' },
    @{ Cell = 'B22'; Value = 'Dai' },
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Cell = 'D22'; Value = '1.00' },
    @{ Cell = 'E22'; Value = '  -0.02%  ' },
    @{ Cell = 'B23'; Value = 'PancakeSwap' },
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' },
    @{ Cell = 'D23'; Value = '2.42' },
    @{ Cell = 'E23'; Value = '  +0.02%  ' },
    @{ Cell = 'B24'; Value = 'Toncoin' },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Cell = 'D24'; Value = '2.28' },
    @{ Cell = 'E24'; Value = '  -5.76%  ' },
    @{ Cell = 'B25'; Value = 'Monero' },
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Cell = 'D25'; Value = '171.12' },
    @{ Cell = 'E25'; Value = '  +1.57%  ' },
    @{ Cell = 'B26'; Value = 'Cosmos' },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Cell = 'D26'; Value = '9.15' },
    @{ Cell = 'E26'; Value = '  +1.12%  ' },
    @{ Cell = 'B27'; Value = 'Kaspa' },
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' },
    @{ Cell = 'D27'; Value = '0.136' },
    @{ Cell = 'E27'; Value = '  -4.63%  ' },
    @{ Cell = 'B28'; Value = 'ImmutableX' },
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Cell = 'D28'; Value = '1.41' },
    @{ Cell = 'E28'; Value = '  -1.67%  ' },
    @{ Cell = 'D29'; Value = '19.50' },
    @{ Cell = 'E29'; Value = '  -0.31%  ' },
    @{ Cell = 'B30'; Value = 'Stellar' },
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Cell = 'D30'; Value = '0.121' },
    @{ Cell = 'E30'; Value = '  +1.38%  ' },
    @{ Cell = 'B31'; Value = 'Filecoin' },
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Cell = 'D31'; Value = '4.71' },
    @{ Cell = 'E31'; Value = '  +0.60%  ' },
    @{ Cell = 'B32'; Value = 'Hedera' },
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Cell = 'D32'; Value = '0.0637' },
    @{ Cell = 'E32'; Value = '  +1.09%  ' },
    @{ Cell = 'B33'; Value = 'InternetComputer(DFINITY)' },
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' },
    @{ Cell = 'D33'; Value = '4.68' },
    @{ Cell = 'E33'; Value = '  +1.06%  ' },
    @{ Cell = 'B34'; Value = 'LidoDAOToken' },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' },
    @{ Cell = 'D34'; Value = '2.48' },
    @{ Cell = 'E34'; Value = '  -5.24%  ' },
    @{ Cell = 'B35'; Value = 'WEMIXToken' },
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'D35'; Value = '1.82' },
    @{ Cell = 'E35'; Value = '  -0.43%  ' },
    @{ Cell = 'B36'; Value = 'RenderToken' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Cell = 'D36'; Value = '3.41' },
    @{ Cell = 'E36'; Value = '  -2.83%  ' },
    @{ Cell = 'B37'; Value = 'BinanceUSD' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' },
    @{ Cell = 'D37'; Value = '1.00' },
    @{ Cell = 'E37'; Value = '  +0.06%  ' },
    @{ Cell = 'B38'; Value = 'THORChain' },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune' },
    @{ Cell = 'D38'; Value = '5.37' },
    @{ Cell = 'E38'; Value = '  -2.56%  ' },
    @{ Cell = 'B39'; Value = 'Cronos' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Cell = 'D39'; Value = '0.0979' },
    @{ Cell = 'E39'; Value = '  -2.03%  ' },
    @{ Cell = 'B40'; Value = 'Aave' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Cell = 'D40'; Value = '99.67' },
    @{ Cell = 'E40'; Value = '  +2.16%  ' },
    @{ Cell = 'B41'; Value = 'HuobiToken' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' },
    @{ Cell = 'D41'; Value = '2.88' },
    @{ Cell = 'E41'; Value = '  -2.32%  ' },
    @{ Cell = 'B42'; Value = 'VeChain' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D42'; Value = '0.0216' },
    @{ Cell = 'E42'; Value = '  -0.17%  ' },
    @{ Cell = 'B43'; Value = 'InjectiveProtocol' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Cell = 'D43'; Value = '16.69' },
    @{ Cell = 'E43'; Value = '  +4.04%  ' },
    @{ Cell = 'B44'; Value = 'Maker' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Cell = 'D44'; Value = '1.442.79' },
    @{ Cell = 'E44'; Value = '  -1.61%  ' },
    @{ Cell = 'B45'; Value = 'TrustWalletToken' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Cell = 'D45'; Value = '1.15' },
    @{ Cell = 'E45'; Value = '  -1.01%  ' },
    @{ Cell = 'B46'; Value = 'FTXToken' },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' },
    @{ Cell = 'D46'; Value = '4.19' },
    @{ Cell = 'E46'; Value = '  +0.95%  ' },
    @{ Cell = 'B47'; Value = 'ARBITRUM' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' },
    @{ Cell = 'D47'; Value = '1.06' },
    @{ Cell = 'E47'; Value = '  -0.91%  ' },
    @{ Cell = 'B48'; Value = 'FraxShare' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Cell = 'D48'; Value = '7.43' },
    @{ Cell = 'E48'; Value = '  +1.36%  ' },
    @{ Cell = 'B49'; Value = 'MXToken' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'D49'; Value = '3.00' },
    @{ Cell = 'E49'; Value = '  -1.43%  ' },
    @{ Cell = 'B50'; Value = 'RocketPoolETH' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' },
    @{ Cell = 'D50'; Value = '2.274.59' },
    @{ Cell = 'E50'; Value = '  -1.97%  ' },
    @{ Cell = 'B51'; Value = 'MultiversX' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld' },
    @{ Cell = 'D51'; Value = '46.63' },
    @{ Cell = 'E51'; Value = '  -0.13%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell -match '^[D]\d+$') {
        $cell.NumberFormat = '@'
    }
    $cell.Value = $u.Value
}
